$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item("ConsorcioDespesas")

$lastExistingRow = 349
$startRow = 350

$data = @(
    ,@(45596, 4, 122, "319011010100", 900.4, 900.4, 900.4)
    ,@(45596, 4, 122, "319013010100", 72.03, 72.03, 73.54)
    ,@(45596, 4, 122, "319013020100", 189.08, 189.08, 193.04)
    ,@(45596, 4, 122, "339014140000", 0, 0, 0)
    ,@(45596, 4, 122, "339030000000", 0, 15.51, 15.51)
    ,@(45596, 4, 122, "339033010000", 0, 0, 0)
    ,@(45596, 4, 122, "339039000000", 47.29, 61.04, 50.26)
    ,@(45596, 4, 122, "339039990100", 59.11, 59.11, 59.11)
    ,@(45596, 4, 122, "339039400000", 0, 111.6, 111.6)
    ,@(45596, 4, 122, "339046010100", 53.66, 53.66, 53.66)
    ,@(45596, 4, 122, "339047000000", 0, 0, 0)
    ,@(45596, 4, 122, "339049010000", 1.96, 1.96, 1.96)
    ,@(45596, 4, 122, "449052000000", 0, 0, 0)
    ,@(45596, 10, 302, "334041390500", 9273.68, 9273.68, 9273.68)
    ,@(45596, 10, 302, "334041391100", 865.54, 865.54, 865.54)
    ,@(45596, 10, 302, "334041391000", 655.91, 655.91, 655.91)
)

$ano = 2024
$r = $startRow
foreach ($row in $data) {
  # Copy number formats (date/currency/code formats) from the last existing row
  $ws.Range("A${lastExistingRow}:K${lastExistingRow}").Copy()
  $ws.Range("A${r}:K${r}").PasteSpecial(-4122)

  $ws.Range("A$r").Value = "COFRON"
  $ws.Range("B$r").Value = $row[0]
  $ws.Range("C$r").Value = $row[1]
  $ws.Range("D$r").Value = $row[2]
  $ws.Range("E$r").Value = [double]$row[3]
  $ws.Range("F$r").Value = $row[4]
  $ws.Range("G$r").Value = $row[5]
  $ws.Range("H$r").Value = $row[6]

  # Calculated columns: set formulas BEFORE the row joins the table
  # (array-entered formulas cannot be assigned to cells already inside a table).
  $ws.Range("I$r").Formula = "=YEAR(ConsorcioDespesas[[#This Row],[data_base]])"
  $ws.Range("J$r").FormulaArray = "=_xlfn.SWITCH(MONTH(ConsorcioDespesas[[#This Row],[data_base]]),1,1,2,1,3,2,4,2,5,3,6,3,7,4,8,4,9,5,10,5,11,6,12,6)"
  $ws.Range("K$r").Formula = "=MONTH(ConsorcioDespesas[[#This Row],[data_base]])"

  $r = $r + 1
}

$endRow = $r - 1
$lo.Resize($ws.Range("A1:K$endRow"))

# Update the view to match where the user scrolled/selected after the edit
$excel.ActiveWindow.ScrollRow = 330
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H366").Select()

$wb.Save()
